$wb = $excel.ActiveWorkbook

# "meta" sheet is the first sheet; "bars-and-line" is the second.
$meta = $wb.Worksheets.Item(1)

# Insert a new row above the existing (empty) row 7 so the previously
# empty A7 (style "1", no value) is pushed down to row 8, and use the
# freed-up row 7 for the new "style" / "default" key-value pair.
$meta.Rows("7:7").Insert()

$meta.Range("A7").Value = "style"
$meta.Range("B7").Value = "default"
